$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the existing header font first (reuses the current style/font slot)
$ws.Range("A1").Font.Bold = $true

# Update header text: A1 becomes "ID", new column B gets "Item Group Name"
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Item Group Name"

# Copy A1's formatting (bold header look) onto B1 so both share the same style
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Resize columns: A narrower, B takes over the old 47-width header look
$ws.Columns.Item(1).ColumnWidth = 35.7109375
$ws.Columns.Item(2).ColumnWidth = 47

# Move the active selection to B3
$ws.Range("B3").Select()
